{"js": "// Corrections des fautes mineures dans le document d'analyse\n// Each entry is [searchText, replacementText]. All searches are done with\n// matchCase:true so we don't accidentally touch casing-similar text, and\n// every search string below is verified unique in the source document.\nconst replacements = [\n  // \"d\u00e9velopper\" -> \"d\u00e9velopp\u00e9\"\n  [\"d\u00e9velopper pour le public\", \"d\u00e9velopp\u00e9 pour le public\"],\n  // remove \"gros \"\n  [\"divis\u00e9 en 4 gros packages\", \"divis\u00e9 en 4 packages\"],\n  // \"doit-\u00eatre\" -> \"doit ... \u00eatre\" + \"voir\" -> \"y acc\u00e9der.\"\n  [\n    \"La page d\\u2019accueil doit-\u00eatre absolument celle du login, afin qu\\u2019aucune personne n\\u2019ayant pas les acc\u00e8s puisse voir \",\n    \"La page d\\u2019accueil doit absolument \u00eatre celle du login, afin qu\\u2019aucune personne n\\u2019ayant pas les acc\u00e8s puisse y acc\u00e9der. \",\n  ],\n  // \"moyens\" -> \"moyen\" + \"avec l'application C#.\" -> \"dans la console d'administration.\"\n  [\n    \"Le seul moyens de les modifier c\\u2019est avec l\\u2019application C#.\",\n    \"Le seul moyen de les modifier c\\u2019est dans la console d\\u2019administration.\",\n  ],\n  // \"info\"/\"infos\" -> \"information\"/\"informations\"\n  [\n    \"Un sous-\u00e9v\u00e8nement ne contient pratiquement aucune info. Ces infos sont en fait un engagement.\",\n    \"Un sous-\u00e9v\u00e8nement ne contient pratiquement aucune information. Ces informations sont en fait un engagement.\",\n  ],\n  // \"envoy\u00e9\" -> \"envoy\u00e9e\" (before \"en fichier XML\")\n  [\"elle doit \u00eatre envoy\u00e9 en fichier XML.\", \"elle doit \u00eatre envoy\u00e9e en fichier XML.\"],\n  // \"devient donc en appel\" -> \"devient donc un appel\"\n  [\"Cette publicit\u00e9 devient donc en appel d\\u2019offre.\", \"Cette publicit\u00e9 devient donc un appel d\\u2019offre.\"],\n  // \"accepter\" -> \"accept\u00e9\"\n  [\"il peut \u00eatre annul\u00e9, accepter, etc\", \"il peut \u00eatre annul\u00e9, accept\u00e9, etc\"],\n  // simplify \"seulement un seul ... qui peut \u00eatre accept\u00e9\" -> \"un seul ... accept\u00e9\"\n  [\n    \"Il est important de noter qu\\u2019il peut y avoir seulement un seul appel d\\u2019offre qui peut \u00eatre accept\u00e9.\",\n    \"Il est important de noter qu\\u2019il peut y avoir un seul appel d\\u2019offre accept\u00e9.\",\n  ],\n  // \"Google maps\" -> \"Google Maps\"\n  [\"Google maps\", \"Google Maps\"],\n  // \"acheter\" -> \"achet\u00e9\"\n  [\"si le billet acheter est en pr\u00e9vente\", \"si le billet achet\u00e9 est en pr\u00e9vente\"],\n  // reorder \"afin de modifier ou lui rappeler de\" -> \"afin de lui rappeler de modifier\"\n  [\n    \"afin de modifier ou lui rappeler de quelque chose.\",\n    \"afin de lui rappeler de modifier quelque chose.\",\n  ],\n  // \"autre\" -> \"utilisateur\"\n  [\n    \"Destinataire\\u00a0: Message \u00e0 envoyer \u00e0 un autre destinataire.\",\n    \"Destinataire\\u00a0: Message \u00e0 envoyer \u00e0 un utilisateur destinataire.\",\n  ],\n  // \"leurs\"->\"le\" and \"ans\"->\"dans\"\n  [\n    \"\u00e9crivent leurs message qu\\u2019ils veulent ans une boite de texte.\",\n    \"\u00e9crivent le message qu\\u2019ils veulent dans une boite de texte.\",\n  ],\n  // comma -> colon, \"supprim\u00e9e\" -> \"supprim\u00e9s\"\n  [\n    \"Les statuts des m\u00e9mos ont 3 formes, lues, non-lues et supprim\u00e9e. \",\n    \"Les statuts des m\u00e9mos ont 3 formes : lues, non-lues et supprim\u00e9s. \",\n  ],\n  // \"statut\" -> \"m\u00e9mo\"\n  [\"Date\\u00a0: Prend la date de la cr\u00e9ation du statut.\", \"Date\\u00a0: Prend la date de la cr\u00e9ation du m\u00e9mo.\"],\n];\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Pattern not found: \" + searchText);\n  }\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// The \"_GoBack\" bookmark (an artifact of Word's last-edit-position marker)\n// moved from right after \"Analyse \" to inside \"Google Maps\" (between the\n// \"M\" and \"aps\"). Reproduce that relocation.\n//\n// Note: calling .delete()/.clear() directly on the empty (collapsed)\n// bookmark range wipes out the whole enclosing paragraph's text in this\n// engine, so instead we remove the bookmark by replacing its paragraph's\n// text with itself (a self-replace drops any bookmark anchors inside it\n// without touching the visible text).\nconst oldBookmark = context.document.getBookmarkRange(\"_GoBack\");\noldBookmark.load(\"isNullObject\");\nawait context.sync();\nif (!oldBookmark.isNullObject) {\n  const bookmarkParagraph = oldBookmark.paragraphs.getFirst();\n  bookmarkParagraph.load(\"text\");\n  await context.sync();\n  bookmarkParagraph.insertText(bookmarkParagraph.text, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst apsRange = context.document.body.search(\"aps\", { matchCase: true });\napsRange.load(\"items\");\nawait context.sync();\nif (apsRange.items.length > 0) {\n  const startOfAps = apsRange.items[0].getRange(\"Start\");\n  startOfAps.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Correction des fautes mineures dans le document d'analyse\n$d = $word.ActiveDocument\n$nbsp = [char]0x00A0\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $result) {\n        throw \"Pattern not found: $old\"\n    }\n}\n\n# \"d\u00e9velopper\" -> \"d\u00e9velopp\u00e9\"\nReplace-Text \"d\u00e9velopper pour le public\" \"d\u00e9velopp\u00e9 pour le public\"\n\n# remove \"gros \"\nReplace-Text \"divis\u00e9 en 4 gros packages\" \"divis\u00e9 en 4 packages\"\n\n# \"doit-\u00eatre\" -> \"doit ... \u00eatre\" + \"voir\" -> \"y acc\u00e9der.\"\nReplace-Text \"La page d\u2019accueil doit-\u00eatre absolument celle du login, afin qu\u2019aucune personne n\u2019ayant pas les acc\u00e8s puisse voir \" \"La page d\u2019accueil doit absolument \u00eatre celle du login, afin qu\u2019aucune personne n\u2019ayant pas les acc\u00e8s puisse y acc\u00e9der. \"\n\n# \"moyens\" -> \"moyen\" + \"avec l'application C#.\" -> \"dans la console d'administration.\"\nReplace-Text \"Le seul moyens de les modifier c\u2019est avec l\u2019application C#.\" \"Le seul moyen de les modifier c\u2019est dans la console d\u2019administration.\"\n\n# \"info\"/\"infos\" -> \"information\"/\"informations\"\nReplace-Text \"Un sous-\u00e9v\u00e8nement ne contient pratiquement aucune info. Ces infos sont en fait un engagement.\" \"Un sous-\u00e9v\u00e8nement ne contient pratiquement aucune information. Ces informations sont en fait un engagement.\"\n\n# \"envoy\u00e9\" -> \"envoy\u00e9e\" (before \"en fichier XML\")\nReplace-Text \"elle doit \u00eatre envoy\u00e9 en fichier XML.\" \"elle doit \u00eatre envoy\u00e9e en fichier XML.\"\n\n# \"devient donc en appel\" -> \"devient donc un appel\"\nReplace-Text \"Cette publicit\u00e9 devient donc en appel d\u2019offre.\" \"Cette publicit\u00e9 devient donc un appel d\u2019offre.\"\n\n# \"accepter\" -> \"accept\u00e9\"\nReplace-Text \"il peut \u00eatre annul\u00e9, accepter, etc\" \"il peut \u00eatre annul\u00e9, accept\u00e9, etc\"\n\n# simplify \"seulement un seul ... qui peut \u00eatre accept\u00e9\" -> \"un seul ... accept\u00e9\"\nReplace-Text \"Il est important de noter qu\u2019il peut y avoir seulement un seul appel d\u2019offre qui peut \u00eatre accept\u00e9.\" \"Il est important de noter qu\u2019il peut y avoir un seul appel d\u2019offre accept\u00e9.\"\n\n# \"Google maps\" -> \"Google Maps\"\nReplace-Text \"Google maps\" \"Google Maps\"\n\n# \"acheter\" -> \"achet\u00e9\"\nReplace-Text \"si le billet acheter est en pr\u00e9vente\" \"si le billet achet\u00e9 est en pr\u00e9vente\"\n\n# reorder \"afin de modifier ou lui rappeler de\" -> \"afin de lui rappeler de modifier\"\nReplace-Text \"afin de modifier ou lui rappeler de quelque chose.\" \"afin de lui rappeler de modifier quelque chose.\"\n\n# \"autre\" -> \"utilisateur\"\nReplace-Text (\"Destinataire\" + $nbsp + \": Message \u00e0 envoyer \u00e0 un autre destinataire.\") (\"Destinataire\" + $nbsp + \": Message \u00e0 envoyer \u00e0 un utilisateur destinataire.\")\n\n# \"leurs\"->\"le\" and \"ans\"->\"dans\"\nReplace-Text \"\u00e9crivent leurs message qu\u2019ils veulent ans une boite de texte.\" \"\u00e9crivent le message qu\u2019ils veulent dans une boite de texte.\"\n\n# comma -> colon, \"supprim\u00e9e\" -> \"supprim\u00e9s\"\nReplace-Text \"Les statuts des m\u00e9mos ont 3 formes, lues, non-lues et supprim\u00e9e. \" \"Les statuts des m\u00e9mos ont 3 formes : lues, non-lues et supprim\u00e9s. \"\n\n# \"statut\" -> \"m\u00e9mo\"\nReplace-Text (\"Date\" + $nbsp + \": Prend la date de la cr\u00e9ation du statut.\") (\"Date\" + $nbsp + \": Prend la date de la cr\u00e9ation du m\u00e9mo.\")\n\n# The \"_GoBack\" bookmark (Word's last-edit-position marker) moved from right\n# after \"Analyse \" to inside \"Google Maps\" (between the \"M\" and \"aps\").\n# Reproduce that relocation: delete the old one, then insert a new one at\n# the new position.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Google Maps\"\n$find2.Execute() | Out-Null\n$mapsRange = $find2.Parent\n$bookmarkStart = $mapsRange.Start + 8  # length of \"Google M\" -> right after the \"M\"\n$bookmarkRange = $d.Range($bookmarkStart, $bookmarkStart)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange) | Out-Null\n\nWrite-Output \"done\"\n"}
